# This script reproduces a data-entry fix in the stock report: for a number
# of duplicate-SKU row groups, the per-row figures (Closing Qty code in col B,
# Purchase/Sale rate in D, MRP/rate in E, Qty in F, Value in G) had been
# shifted by one row relative to the correct record. The fix rotates the
# B/D/E/F/G values within each affected group of rows (which share the same
# item description in column C) up by one row, wrapping the first row's
# values around to the last row of the group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is a contiguous block of worksheet rows (sharing the same
# item in column C) whose B/D/E/F/G values need to be cyclically rotated.
$groups = @(
    @(149, 150),
    @(264, 265),
    @(279, 280),
    @(316, 317, 318),
    @(372, 373),
    @(379, 380),
    @(382, 383),
    @(421, 422),
    @(536, 537),
    @(590, 591),
    @(599, 600),
    @(601, 602),
    @(709, 710),
    @(720, 721),
    @(859, 860)
)

$cols = @("B", "D", "E", "F", "G")

foreach ($group in $groups) {
    $n = $group.Length

    # Snapshot the current values for the columns we are about to rotate.
    # (Value2 is used for reading — it returns the raw scalar — as opposed
    # to Value, which is unreliable for reads in this host.)
    $snapshot = @{}
    foreach ($row in $group) {
        $rowValues = @{}
        foreach ($col in $cols) {
            $rowValues[$col] = $ws.Range("$col$row").Value2
        }
        $snapshot[$row] = $rowValues
    }

    # Row i receives the values that originally belonged to row i+1 (wrapping
    # around to the first row of the group for the last row).
    for ($i = 0; $i -lt $n; $i++) {
        $targetRow = $group[$i]
        $sourceRow = $group[($i + 1) % $n]
        foreach ($col in $cols) {
            $ws.Range("$col$targetRow").Value2 = $snapshot[$sourceRow][$col]
        }
    }
}
